$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.284.98"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.491.52"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.20"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.23"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.084.80"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.488.64"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.74"
$ws.Range("E16").Value = "  -7.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.335.82"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.62"
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.09"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.628.56"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.78"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("E31").Value = "  -5.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.23"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.510.02"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.39"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.13"
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.89"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.13"
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0779"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.18"
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.38"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.455.40"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.75"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.891"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("E51").Value = "  -1.51%  "
